$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01168537388324469
$ws.Range("C2").Value = 0.05127383665300928

$ws.Range("B3").Value = 0.07226372149614482
$ws.Range("C3").Value = 0.09497099886898824

$ws.Range("B4").Value = 0.08210775098986554
$ws.Range("C4").Value = 0.0564290864261125

$ws.Range("B5").Value = 0.6391606537483809
$ws.Range("C5").Value = 0.5008755618547541

$ws.Range("B6").Value = 0.9542407599495295
$ws.Range("C6").Value = 0.906915053152735

$ws.Range("B7").Value = 0.9496574571092881
$ws.Range("C7").Value = 0.7969969073460706

$ws.Range("B8").Value = 0.7119921400426712
$ws.Range("C8").Value = 0.4999629988120595

$ws.Range("B9").Value = 0.003900811076164246
$ws.Range("C9").Value = 0.03246632099151611

$ws.Range("B10").Value = 0.8478399931966736
$ws.Range("C10").Value = 0.763430302268425
